$wb = $excel.ActiveWorkbook

# New identifiers introduced by this commit (a fresh handoff + a failed-transform entry).
$newMdFile      = "cf29bef5-f826-4efd-9ef1-ae17f2d7cb29.md"
$failedMdFile   = "792e9cb1-48e5-427b-ac7b-7e55e17e34a2.md"
$configFile     = ".localization-config"

$zhXlfFile = "cf29bef5-f826-4efd-9ef1-ae17f2d7cb29.7e9ccfd4c0ff5679616a2a4bb833aaca5f1430f6.zh-cn.xlf"
$deXlfFile = "cf29bef5-f826-4efd-9ef1-ae17f2d7cb29.7e9ccfd4c0ff5679616a2a4bb833aaca5f1430f6.de-de.xlf"

$zhHandoffTime = "2016-01-13 11:32:51"
$deHandoffTime = "2016-01-13 11:33:05"

$readyStatus  = "Ready for handoff"
$failedStatus = "Handoff transform failed"
$notLocalized = "Not to be localized"
$epoch        = "0001-01-01 00:00:00"
$includeTxt   = "Include"
$ignoredTxt   = "Ignored"

$baseRepoUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/bcc48e979129b92dd6f0cf50ac1083962b4818a5"
$zhHandoffRepo = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb5c2882e844af88230a88edff99b14600d22232/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho"
$deHandoffRepo = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e666710a49171cc8180384b5d7d5467a9d788d62/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho"

# ============================================================
# Sheet "Overview": File Name | zh-cn | de-de
# ============================================================
$ws = $wb.Worksheets.Item("Overview")

# Drop every hyperlink on the sheet so we can rebuild them cleanly at their
# final addresses (the engine's hyperlink-delete only works sheet-wide).
$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("A2").Value = $newMdFile
$ws.Range("B2").Value = $readyStatus
$ws.Range("C2").Value = $readyStatus

$ws.Range("A3").Value = $failedMdFile
$ws.Range("B3").Value = $failedStatus
$ws.Range("C3").Value = $failedStatus

$ws.Range("A4").Value = $configFile
$ws.Range("B4").Value = $notLocalized
$ws.Range("C4").Value = $notLocalized

$ws.Hyperlinks.Add($ws.Range("A2"), "$baseRepoUrl/e2e/$newMdFile", [Type]::Missing, [Type]::Missing, $newMdFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$baseRepoUrl/e2e/$failedMdFile", [Type]::Missing, [Type]::Missing, $failedMdFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "$baseRepoUrl/$configFile", [Type]::Missing, [Type]::Missing, $configFile) | Out-Null

# ============================================================
# Sheet "zh-cn"
# ============================================================
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("A2").Value = $newMdFile
$ws.Range("B2").Value = $readyStatus
$ws.Range("C2").Value = $zhXlfFile
$ws.Range("D2").Value = $zhHandoffTime
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = $includeTxt

$ws.Range("A3").Value = $failedMdFile
$ws.Range("B3").Value = $failedStatus
$ws.Range("D3").Value = $epoch
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = $ignoredTxt

$ws.Range("A4").Value = $configFile
$ws.Range("B4").Value = $notLocalized
$ws.Range("D4").Value = $epoch
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = $ignoredTxt

$ws.Hyperlinks.Add($ws.Range("A2"), "$baseRepoUrl/e2e/$newMdFile", [Type]::Missing, [Type]::Missing, $newMdFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "$zhHandoffRepo/$zhXlfFile", [Type]::Missing, [Type]::Missing, $zhXlfFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$baseRepoUrl/e2e/$failedMdFile", [Type]::Missing, [Type]::Missing, $failedMdFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "$baseRepoUrl/$configFile", [Type]::Missing, [Type]::Missing, $configFile) | Out-Null

# ============================================================
# Sheet "de-de"
# ============================================================
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("A2").Value = $newMdFile
$ws.Range("B2").Value = $readyStatus
$ws.Range("C2").Value = $deXlfFile
$ws.Range("D2").Value = $deHandoffTime
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = $includeTxt

$ws.Range("A3").Value = $failedMdFile
$ws.Range("B3").Value = $failedStatus
$ws.Range("D3").Value = $epoch
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = $ignoredTxt

$ws.Range("A4").Value = $configFile
$ws.Range("B4").Value = $notLocalized
$ws.Range("D4").Value = $epoch
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = $ignoredTxt

$ws.Hyperlinks.Add($ws.Range("A2"), "$baseRepoUrl/e2e/$newMdFile", [Type]::Missing, [Type]::Missing, $newMdFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "$deHandoffRepo/$deXlfFile", [Type]::Missing, [Type]::Missing, $deXlfFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$baseRepoUrl/e2e/$failedMdFile", [Type]::Missing, [Type]::Missing, $failedMdFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "$baseRepoUrl/$configFile", [Type]::Missing, [Type]::Missing, $configFile) | Out-Null

Write-Host "Localization status report updated: added 'Handoff transform failed' row and refreshed handoff ids."
